# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in A1 (05:22 -> 05:52)
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 05:52"

# Update Paraguay row (row 124): active cases, critical cases, deaths
$ws.Range("E124").Value = 96
$ws.Range("G124").Value = 2
$ws.Range("H124").Value = 5

# Re-sorted block of small countries (rows 170-176): new totals pushed
# Mongolia and Fiyi up the table, displacing Dominica/Santa Lucia/Angola/
# Liberia/Granada down by one row each.
$ws.Range("A170").Value = "Mongolia"
$ws.Range("B170").Value = 15
$ws.Range("C170").Value = 1
$ws.Range("D170").Value = 2
$ws.Range("E170").Value = 13
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 0

$ws.Range("A171").Value = "Fiyi"
$ws.Range("B171").Value = 14
$ws.Range("C171").Value = 2
$ws.Range("D171").Value = 0
$ws.Range("E171").Value = 14
$ws.Range("F171").Value = 0
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = 0

$ws.Range("A172").Value = "Dominica"
$ws.Range("B172").Value = 14
$ws.Range("C172").Value = 0
$ws.Range("D172").Value = 0
$ws.Range("E172").Value = 14
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 0

$ws.Range("A173").Value = "Santa Lucia"
$ws.Range("B173").Value = 14
$ws.Range("C173").Value = 0
$ws.Range("D173").Value = 1
$ws.Range("E173").Value = 13
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 0

$ws.Range("A174").Value = "Angola"
$ws.Range("B174").Value = 14
$ws.Range("C174").Value = 0
$ws.Range("D174").Value = 2
$ws.Range("E174").Value = 10
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 2

$ws.Range("A175").Value = "Liberia"
$ws.Range("B175").Value = 13
$ws.Range("C175").Value = 0
$ws.Range("D175").Value = 3
$ws.Range("E175").Value = 7
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 3

$ws.Range("A176").Value = "Granada"
$ws.Range("B176").Value = 12
$ws.Range("C176").Value = 0
$ws.Range("D176").Value = 0
$ws.Range("E176").Value = 12
$ws.Range("F176").Value = 2
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 0
